$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 was previously blank; fill it in with a new "Current Fiscal Year" field,
# matching the look (font/alignment) of the other question rows (A7:A9).
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10").Value = "Current Fiscal Year"
$ws.Range("B10").Value = 2023
